{"js": "// The author fixed a typo (\"No}\" -> \"No]\") in the closing instructional\n// paragraph. Word silently tracks the position of the most recent edit\n// with a hidden bookmark named \"_GoBack\"; since this edit is now the\n// most recent one, that bookmark moves from its old location (the start\n// of the \"Sample Output - Gui:\" heading) to right after the fixed text.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// Drop the stale \"_GoBack\" bookmark first - Word only ever keeps a single\n// bookmark with a given name, and it will be re-inserted at the new spot\n// below (matches the diff removing the old bookmarkStart/bookmarkEnd pair).\ndoc.deleteBookmark(\"_GoBack\");\n\n// Locate the full run of text that needs the typo fixed (include the\n// leading tab character so the whole original run is replaced cleanly).\nconst target = body.search(\n  \"\\tSelect [No} to finalize the orders. User may select [Yes] to repeat the loop and commit arbitrarily more orders.\",\n  { matchCase: true }\n);\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  const run = target.items[0];\n\n  // Replace it with the corrected text, split into two runs with the\n  // relocated \"_GoBack\" bookmark sitting between them - exactly mirroring\n  // the target markup: \"...<w:tab/><w:t>Select [No]</w:t></w:r>\n  // <w:bookmarkStart.../><w:bookmarkEnd.../><w:r><w:t> to finalize...\".\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" +\n    '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/><w:t>Select [No]</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t xml:space=\"preserve\"> to finalize the orders. User may select [Yes] to repeat the loop and commit arbitrarily more orders.</w:t></w:r>' +\n    \"</w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\";\n\n  run.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The author fixed a typo (\"No}\" -> \"No]\") in the closing instructional\n# paragraph. Word silently tracks the position of the most recent edit\n# with a hidden bookmark named \"_GoBack\"; since this edit is now the\n# most recent one, that bookmark moves from its old location (the start\n# of the \"Sample Output - Gui:\" heading) to right after the fixed text.\n\n$d = $word.ActiveDocument\n\n# Drop the stale \"_GoBack\" bookmark first - Word only ever keeps a single\n# bookmark with a given name, and it will be re-inserted at the new spot\n# below (matches the diff removing the old bookmarkStart/bookmarkEnd pair).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Locate the full run of text that needs the typo fixed (include the\n# leading tab character so the whole original run is replaced cleanly).\n$r = $d.Content\n$found = $r.Find.Execute(\"`tSelect [No} to finalize the orders. User may select [Yes] to repeat the loop and commit arbitrarily more orders.\")\n\nif ($found) {\n    $target = $d.Range($r.Start, $r.End)\n\n    # Replace it with the corrected text, split into two runs with the\n    # relocated \"_GoBack\" bookmark sitting between them - exactly mirroring\n    # the target markup: \"...<w:tab/><w:t>Select [No]</w:t></w:r>\n    # <w:bookmarkStart.../><w:bookmarkEnd.../><w:r><w:t> to finalize...\".\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p>' +\n        '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/><w:t>Select [No]</w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n        '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t xml:space=\"preserve\"> to finalize the orders. User may select [Yes] to repeat the loop and commit arbitrarily more orders.</w:t></w:r>' +\n        '</w:p></w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $target.InsertXML($xml)\n}\n"}
